# Append 7 new daily COVID summary rows (2022-02-14 .. 2022-02-20) to the
# "overview"/United Kingdom series, continuing directly after the existing
# last row (551, dated 2022-02-13). Extends the sheet's used range from
# A1:H551 to A1:H558.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# date, areaType, areaCode, areaName, cumCasesByPublishDate,
# newCasesByPublishDate, newDeaths28DaysByPublishDate, cumDeaths28DaysByPublishDate
$newRows = @(
    @("2022-02-14","overview","K02000001","United Kingdom",18348029,41648, 35,159605),
    @("2022-02-15","overview","K02000001","United Kingdom",18393951,46186,234,159839),
    @("2022-02-16","overview","K02000001","United Kingdom",18447706,54218,199,160038),
    @("2022-02-17","overview","K02000001","United Kingdom",18499058,51899,183,160221),
    @("2022-02-18","overview","K02000001","United Kingdom",18546205,47685,158,160379),
    @("2022-02-19","overview","K02000001","United Kingdom",18580216,34377,128,160507),
    @("2022-02-20","overview","K02000001","United Kingdom",18605752,25696, 74,160581)
)

$startRow = 552
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    # Column A holds an ISO date string ("2022-02-14"), not a real Excel
    # date. A leading apostrophe forces text entry (same as typing it into
    # Excel); resetting the style back to "Normal" afterwards drops the
    # quote-prefix formatting so the cell ends up with the default style,
    # matching the rest of the column.
    $ws.Range("A$r").Value = "'" + $row[0]
    $ws.Range("A$r").Style = "Normal"

    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
    $ws.Range("F$r").Value = $row[5]
    $ws.Range("G$r").Value = $row[6]
    $ws.Range("H$r").Value = $row[7]
}
